$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (rows 2-398) holds a "Förändrad" (changed) date that was bumped
# by one day, from 2023-09-02 (45171) to 2023-09-03 (45172).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 398 }

$range = $ws.Range("C2:C$lastRow")
foreach ($cell in $range.Cells) {
    if ($cell.Value2 -eq 45171) {
        $cell.Value2 = 45172
    }
}
